# "cloud integration has been done" -- refresh the generated test-data
# value in the "Test Data" sheet (dataSheet/getDataExcel.xlsx) with the
# newly pulled FireFlink project identifier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")

$ws.Range("B2").Value = "FireFlink_26884"
